# Transportation Technology Logit Exponents.xlsx
# EPS v3.3.1 -> v3.4.2 update
#
# Changes:
#  - "About" sheet: rewritten Notes/Sources block (old EPA citation replaced
#    with a shorter explanation of the -5 / -8 / -10 calibration values).
#  - "TTLE" sheet: Freight column for LDVs changed from -5 to -8, and for
#    HDVs from -5 to -10 (all other logit exponents stay -5); header/data
#    cell formatting simplified (no more per-cell overrides on the data
#    block).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")

# Wipe the sheet completely (content + formatting) and rebuild it from
# scratch so no stale rows/styles are left behind.
$about.Rows("1:69").Delete() | Out-Null

$about.Range("A1").Value = "TTLE Transportation Technology Logit Exponents"
$about.Range("A1").Font.Bold = $true

$about.Range("A3").Value = "Sources:"
$about.Range("A3").Font.Bold = $true
$about.Range("B3").Value = "Calibration"

$about.Range("A5").Value = "Notes"
$about.Range("A5").Font.Bold = $true

$about.Range("A6").Value = "The logit exponents express how large of a cost difference between technology options"
$about.Range("A7").Value = "is required to produce a change in technology selection.  This parameter needs to be"
$about.Range("A8").Value = "obtained via model calibration - e.g. testing a given price intervention with different"
$about.Range("A9").Value = "logit exponent values until it produces a technology choice shift that matches real-world"
$about.Range("A10").Value = "data on technology buyers' behavior."

$about.Range("A12").Value = "We choose a value of -5 for most vehicle types except freight LDVs and HDVs"
$about.Range("A13").Value = "which we assign values of -8 and -10 due to their larger price sensitivity."

$about.Range("A15").Value = "For more on this, see the ""Modified Logit"" equation description at:"
$about.Range("A16").Value = "https://jgcri.github.io/gcam-doc/choice.html"

# Recreate the lone formatted-but-empty cell at the bottom of the sheet
# (keeps the sheet dimension at A1:B54, matching the source file).
$about.Range("A54").Font.Bold = $true

# ---------------------------------------------------------------------
# "TTLE" sheet
# ---------------------------------------------------------------------
$ttle = $wb.Worksheets.Item("TTLE")

# Freight logit exponents: LDVs -5 -> -8, HDVs -5 -> -10
$ttle.Range("C2").Value = -8
$ttle.Range("C3").Value = -10

# Remove the per-cell style overrides on the data block (B2:C7) so the
# cells fall back to the default (general) style, matching the cleaned
# up workbook.
$ttle.Range("B2:C7").ClearFormats() | Out-Null

# Re-normalize the header row styling (duplicate cellXfs collapsed during
# cleanup): A1 stays bold + wrap text, B1/C1 stay right-aligned.
$ttle.Range("A1").ClearFormats() | Out-Null
$ttle.Range("A1").Font.Bold = $true
$ttle.Range("A1").WrapText = $true

$ttle.Range("B1:C1").ClearFormats() | Out-Null
$ttle.Range("B1:C1").HorizontalAlignment = -4152  # xlRight

$ttle.Range("C4").Select() | Out-Null

# Leave "About" as the active sheet/tab (selected last so it stays active).
$about.Activate() | Out-Null
$about.Range("A14").Select() | Out-Null
